$d = $word.ActiveDocument

# --- 1. Insert new paragraph "Rename 'Resources' menu item to 'Admin' and
#        move 'Utilities' menu item to admin" right after the paragraph
#        "Send emails asynchronously" (i.e. right before the "Medium"
#        paragraph). Inserting after that run copies its paragraph / run
#        formatting (ilvl=1, numId=3, HTMLCode rStyle, minor theme fonts,
#        sz/szCs 22) which matches the target markup exactly.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Send emails asynchronously*") {
        $anchor = $p
    }
}

if ($anchor -ne $null) {
    $anchorIndex = $anchor.Index
    $anchor.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($anchorIndex + 1)
    $newPara.Range.Text = "Rename ‘Resources’ menu item to ‘Admin’ and move ‘Utilities’ menu item to admin"
}

# --- 2. Remove the paragraph "Add change owner utility function" (it used
#        to sit right after "Disable edit button in appropriate expense
#        view screen").
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Add change owner utility function*") {
        $p.Range.Delete()
    }
}
